# Lab 3: add a new empty paragraph right after the displayed formula
#   Var(Y_ij) = sigma^2_0B + year^2_ij x sigma^2_1B + 2 x year_ij x sigma_01 + sigma^2_W
# (the oMathPara that follows "Use the formula below to compute the estimate of ...").
# The new paragraph carries the same indent (left=720, firstLine=720) and run
# formatting (Times New Roman, szCs=24) as the formula paragraph itself -- i.e.
# exactly what Word produces when you put the cursor at the end of that
# paragraph and press Enter.

$d = $word.ActiveDocument

# Anchor on the unique sentence that immediately precedes the formula, then
# walk forward: [anchor paragraph] -> (blank paragraph) -> (formula paragraph).
$anchorRange = $d.Content
$anchorRange.Find.Execute(
    "Use the formula below to compute the estimate of",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$anchorParagraph = $anchorRange.Paragraphs(1)
$formulaParagraph = $anchorParagraph.Next().Next()

# Insert a brand-new (no prior revision ids) empty paragraph immediately after
# the formula paragraph, without disturbing the formula paragraph itself.
$insertionPoint = $d.Range($formulaParagraph.Range.End, $formulaParagraph.Range.End)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr>' +
        '<w:ind w:left="720" w:firstLine="720"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>' +
            '<w:szCs w:val="24"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
'</w:p>'

$insertionPoint.InsertXML($newParagraphXml) | Out-Null
